$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 7
$ws.Range("D7").Value = "LangChain이란?"
$ws.Range("E7").Value = "https://jayhey.github.io/deep%20learning/2023/04/23/langchain/"

# Row 9
$ws.Range("D9").Value = "[SIAI공지] 장학 지원 관련 옵션"
$ws.Range("E9").Value = "https://pdsi.pabii.com/siai-funding-options-2023/#utm_source=rss&utm_medium=rss&utm_campaign=siai-funding-options-2023"

# Row 28
$ws.Range("D28").Value = "[논문 리뷰]Towards a Natural Motion Generator: a Pipeline to Control a Humanoid based on Motion Data"
$ws.Range("E28").Value = "https://ropiens.tistory.com/217"

# Row 32
$ws.Range("D32").Value = "[GIT] pre-commit을 이용한 코드 스타일 관리"
$ws.Range("E32").Value = "https://dodonam.tistory.com/426"

$wb.Save()
